# Vaatimusmäärityksen, muistiinpanojen ja tuntipäiväkirjan päivittäminen -KV
#
# Fills in Katja's (R/S/T columns) time-tracking block for rows 14-17 with
# the newly logged work sessions, then restores the sheet's on-screen
# selection to match where the author ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New log entries for Katja (columns R=Päivämäärä, S=Käytetty aika
#     tunneissa, T=Selite). Shared-string text is written first, in the
#     order the new entries appear so the newly created shared-string
#     table entries line up (37: "Vaatimusmäärittelyn tekemistä",
#     38: "...ja Sprintin kirjaamista", 39: "Scrumiin tutustuminen...").

$ws.Range("T15").Value = "Vaatimusmäärittelyn tekemistä"
$ws.Range("T16").Value = "Vaatimusmäärittelyn tekemistä"
$ws.Range("T17").Value = "Vaatimusmäärittelyn tekemistä ja Sprintin kirjaamista"
$ws.Range("T14").Value = "Scrumiin tutustuminen, käsikirjan ja vaatimusmäärittelyn lukemista"

# Dates (stored as serial numbers, formatted by the existing cell style)
$ws.Range("R14").Value = 45340
$ws.Range("R15").Value = 45342
$ws.Range("R16").Value = 45344
$ws.Range("R17").Value = 45345

# Hours used
$ws.Range("S14").Value = 2
$ws.Range("S15").Value = 1
$ws.Range("S16").Value = 4
$ws.Range("S17").Value = 4

# --- Restore on-screen selection/view state ---
$ws.Activate()
[void]$ws.Range("N43").Select()
